# Auto-generated edit script: updates profit-calc values on the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
# to match the refreshed Universalis market-board snapshot pulled in by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 5497.4287
$ws.Range("I32").Value = 4430
$ws.Range("J32").Value = 5788.5454
$ws.Range("K32").Value = 4430
$ws.Range("L32").Value = 5788.5454
$ws.Range("M32").Value = -4104
$ws.Range("N32").Value = -6440.5454
$ws.Range("H125").Value = 5399.8335
$ws.Range("I125").Value = 3799.6667
$ws.Range("J125").Value = 7000
$ws.Range("K125").Value = 34197.0003
$ws.Range("L125").Value = 63000
$ws.Range("M125").Value = -31737.0003
$ws.Range("N125").Value = -67920

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1056
$ws.Range("I2").Value = 905.875
$ws.Range("J2").Value = 1656.5
$ws.Range("K2").Value = 905.875
$ws.Range("L2").Value = 1656.5
$ws.Range("M2").Value = -792.875
$ws.Range("N2").Value = -1882.5
$ws.Range("H63").Value = 4993.636
$ws.Range("I63").Value = 2682.1538
$ws.Range("J63").Value = 8332.444
$ws.Range("K63").Value = 2682.1538
$ws.Range("L63").Value = 8332.444
$ws.Range("M63").Value = -1996.1538
$ws.Range("N63").Value = -9704.444
$ws.Range("H66").Value = 4993.636
$ws.Range("I66").Value = 2682.1538
$ws.Range("J66").Value = 8332.444
$ws.Range("K66").Value = 13410.769
$ws.Range("L66").Value = 41662.22
$ws.Range("M66").Value = -9978.769
$ws.Range("N66").Value = -48526.22
$ws.Range("H114").Value = 70000
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 70000
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 70000
$ws.Range("N114").Value = -78678
$ws.Range("H116").Value = 1056
$ws.Range("I116").Value = 905.875
$ws.Range("J116").Value = 1656.5
$ws.Range("K116").Value = 905.875
$ws.Range("L116").Value = 1656.5
$ws.Range("M116").Value = 1388.125
$ws.Range("N116").Value = -6244.5
$ws.Range("H132").Value = 41669736
$ws.Range("I132").Value = 2483.2354
$ws.Range("J132").Value = 142861650
$ws.Range("K132").Value = 7449.706200000001
$ws.Range("L132").Value = 428584950
$ws.Range("M132").Value = -4919.706200000001
$ws.Range("N132").Value = -428590010

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1056
$ws.Range("I3").Value = 905.875
$ws.Range("J3").Value = 1656.5
$ws.Range("K3").Value = 905.875
$ws.Range("L3").Value = 1656.5
$ws.Range("M3").Value = -791.875
$ws.Range("N3").Value = -1884.5
$ws.Range("H86").Value = 16127.637
$ws.Range("I86").Value = 16200
$ws.Range("J86").Value = 16093.866
$ws.Range("K86").Value = 16200
$ws.Range("L86").Value = 16093.866
$ws.Range("M86").Value = -15077
$ws.Range("N86").Value = -18339.866
$ws.Range("H89").Value = 16127.637
$ws.Range("I89").Value = 16200
$ws.Range("J89").Value = 16093.866
$ws.Range("K89").Value = 81000
$ws.Range("L89").Value = 80469.33
$ws.Range("M89").Value = -75384
$ws.Range("N89").Value = -91701.33
$ws.Range("H105").Value = 8088.4375
$ws.Range("I105").Value = 12787.667
$ws.Range("J105").Value = 2046.5714
$ws.Range("K105").Value = 12787.667
$ws.Range("L105").Value = 2046.5714
$ws.Range("M105").Value = -11040.667
$ws.Range("N105").Value = -5540.5714
$ws.Range("H134").Value = 2973.8572
$ws.Range("I134").Value = 2665.7666
$ws.Range("J134").Value = 4822.4
$ws.Range("K134").Value = 7997.2998
$ws.Range("L134").Value = 14467.2
$ws.Range("M134").Value = -5462.2998
$ws.Range("N134").Value = -19537.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 47626270
$ws.Range("I31").Value = 5295
$ws.Range("J31").Value = 142868220
$ws.Range("K31").Value = 5295
$ws.Range("L31").Value = 142868220
$ws.Range("M31").Value = -5000
$ws.Range("N31").Value = -142868810
$ws.Range("H34").Value = 47626270
$ws.Range("I34").Value = 5295
$ws.Range("J34").Value = 142868220
$ws.Range("K34").Value = 5295
$ws.Range("L34").Value = 142868220
$ws.Range("M34").Value = -5093
$ws.Range("N34").Value = -142868624
$ws.Range("H134").Value = 1202.7778
$ws.Range("I134").Value = 1220.8823
$ws.Range("J134").Value = 895
$ws.Range("K134").Value = 3662.6469
$ws.Range("L134").Value = 2685
$ws.Range("M134").Value = -1127.6469
$ws.Range("N134").Value = -7755

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2419.0833
$ws.Range("I113").Value = 1949.5
$ws.Range("J113").Value = 2888.6667
$ws.Range("K113").Value = 5848.5
$ws.Range("L113").Value = 8666.000100000001
$ws.Range("M113").Value = -3678.5
$ws.Range("N113").Value = -13006.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4146.0586
$ws.Range("I70").Value = 3963.3572
$ws.Range("J70").Value = 4998.6665
$ws.Range("K70").Value = 3963.3572
$ws.Range("L70").Value = 4998.6665
$ws.Range("M70").Value = -3693.3572
$ws.Range("N70").Value = -5538.6665
$ws.Range("H73").Value = 4146.0586
$ws.Range("I73").Value = 3963.3572
$ws.Range("J73").Value = 4998.6665
$ws.Range("K73").Value = 3963.3572
$ws.Range("L73").Value = 4998.6665
$ws.Range("M73").Value = -3027.3572
$ws.Range("N73").Value = -6870.6665
$ws.Range("H102").Value = 1982.6061
$ws.Range("I102").Value = 1463.2
$ws.Range("J102").Value = 2781.6924
$ws.Range("K102").Value = 1463.2
$ws.Range("L102").Value = 2781.6924
$ws.Range("M102").Value = 158.8
$ws.Range("N102").Value = -6025.6924
$ws.Range("H139").Value = 109982.336
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 109982.336
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 109982.336
$ws.Range("N139").Value = -120262.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1710.9269
$ws.Range("I46").Value = 619.29034
$ws.Range("J46").Value = 5095
$ws.Range("K46").Value = 619.29034
$ws.Range("L46").Value = 5095
$ws.Range("M46").Value = -431.29034
$ws.Range("N46").Value = -5471
$ws.Range("H104").Value = 25000
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 25000
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 25000
$ws.Range("N104").Value = -31988
$ws.Range("H119").Value = 85000
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 85000
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 85000
$ws.Range("N119").Value = -94676
$ws.Range("H136").Value = 2346.8235
$ws.Range("I136").Value = 1758.5172
$ws.Range("J136").Value = 5759
$ws.Range("K136").Value = 5275.5516
$ws.Range("L136").Value = 17277
$ws.Range("M136").Value = -2725.5516
$ws.Range("N136").Value = -22377

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("H122").Value = 31251682
$ws.Range("I122").Value = 40001304
$ws.Range("J122").Value = 3028.8572
$ws.Range("K122").Value = 120003912
$ws.Range("L122").Value = 9086.571599999999
$ws.Range("M122").Value = -120001462
$ws.Range("N122").Value = -13986.5716
$ws.Range("H132").Value = 4618.778
$ws.Range("I132").Value = 4736.143
$ws.Range("J132").Value = 4208
$ws.Range("K132").Value = 14208.429
$ws.Range("L132").Value = 12624
$ws.Range("M132").Value = -11678.429
$ws.Range("N132").Value = -17684
$ws.Range("H136").Value = 2978.7856
$ws.Range("I136").Value = 2217.6667
$ws.Range("J136").Value = 4348.8
$ws.Range("K136").Value = 6653.000100000001
$ws.Range("L136").Value = 13046.4
$ws.Range("M136").Value = -4103.000100000001
$ws.Range("N136").Value = -18146.4

# Rows where the HQ-profit figure no longer applies (NQ-only Leve) -- drop the stale N column value
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N63").ClearContents()
$ws.Range("N66").ClearContents()
